$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# Walk the paragraphs once and patch the ones that need structural
# changes, rebuilding each as a minimal, equivalent <w:p> fragment via
# Range.InsertXML so the existing run/formatting layout (bold labels,
# italic lead-ins, xml:space="preserve" text) is kept intact.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text

    if ($text -eq "Write a 3 page report on a mathematical topic. Your report is to be written in LaTeX and must use aspects of programming (Python and/or Sage) to illustrate the particular topic.`r") {
        # Add a sentence describing the intended audience of the report.
        $frag = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">Write a 3 page report on a mathematical topic. Consider the target audience of your report to be first year mathematics students wanting to learn about a given topic. Your report is to be written in LaTeX and must use aspects of programming (Python and/or Sage) to illustrate the particular topic.</w:t></w:r></w:p>"
        $para.Range.InsertXML($frag)
    }
    elseif ($text -like "You are encouraged to choose a topic,*") {
        # "choose a topic" -> "choose your own topic" (keep the italic lead-in
        # run and the plain continuation run untouched otherwise).
        $frag = "<w:p $wNs><w:r><w:rPr><w:i/></w:rPr><w:t xml:space=`"preserve`">You are encouraged to choose your own topic</w:t></w:r><w:r><w:t xml:space=`"preserve`">, if you do so I recommend checking with me (Vince Knight) that the topic is appropriate. If you are unable to choose a topic select one from the following:</w:t></w:r></w:p>"
        $para.Range.InsertXML($frag)
    }
    elseif ($text -like "Code*(50% weighting)*") {
        # Give the "Code" bullet its own numbering id (4), matching the id
        # already used by its own sub-bullets instead of sharing id 3 with
        # the unrelated "Presentation" bullet.
        $frag = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"4`"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`">Code</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t xml:space=`"preserve`">(50% weighting)</w:t></w:r></w:p>"
        $para.Range.InsertXML($frag)
    }
    elseif ($text -like "Content*(30% weighting)*") {
        # Same idea: give "Content" its own numbering id (5).
        $frag = "<w:p $wNs><w:pPr><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"5`"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=`"preserve`">Content</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t xml:space=`"preserve`">(30% weighting)</w:t></w:r></w:p>"
        $para.Range.InsertXML($frag)
    }
}

# ---------------------------------------------------------------------
# Rename the custom styles so their display names match current Word
# naming conventions (the underlying style ids are immutable via the
# object model - only the display/local name can be changed).
# ---------------------------------------------------------------------
$d.Styles("ImageCaption").NameLocal = "Picture Caption"
$d.Styles("FootnoteRef").NameLocal = "Footnote Reference"
$d.Styles("Link").NameLocal = "Hyperlink"
